# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period table (rows 16-23) used to be grouped by worker
# (each worker's "1802" row followed by their "1801" row). The update
# regroups it by period instead: all 4 workers for period "1801" first,
# then the same 4 workers for period "1802".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 16; Doc = "1143363072"; Name = "CINDY DAYANA MADRID ORTIZ";       Periodo = "1801"; Mora = 80000;  Salario = 2000000 },
    @{ Row = 17; Doc = "1129523066"; Name = "BRENDA SOFIA STRUX MONTERROZA";   Periodo = "1801"; Mora = 48000;  Salario = 1200000 },
    @{ Row = 18; Doc = "1143332454"; Name = "PABLO RAFAEL HERRERA CAPDEVILLA"; Periodo = "1801"; Mora = 100000; Salario = 2500000 },
    @{ Row = 19; Doc = "1047435092"; Name = "ANA ISABEL VELASCO BARRETO";      Periodo = "1801"; Mora = 80000;  Salario = 2000000 },
    @{ Row = 20; Doc = "1143363072"; Name = "CINDY DAYANA MADRID ORTIZ";       Periodo = "1802"; Mora = 80000;  Salario = 2000000 },
    @{ Row = 21; Doc = "1129523066"; Name = "BRENDA SOFIA STRUX MONTERROZA";   Periodo = "1802"; Mora = 48000;  Salario = 1200000 },
    @{ Row = 22; Doc = "1143332454"; Name = "PABLO RAFAEL HERRERA CAPDEVILLA"; Periodo = "1802"; Mora = 100000; Salario = 2500000 },
    @{ Row = 23; Doc = "1047435092"; Name = "ANA ISABEL VELASCO BARRETO";      Periodo = "1802"; Mora = 80000;  Salario = 2000000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Name
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Mora
    $ws.Cells.Item($r, 7).Value = $item.Salario
}

# Columns B:J are set to "best fit" widths; since the text in columns
# C:G changed length, refresh their widths to fit the new content.
$widths = @{
    2  = 16.072916666666668
    3  = 9.983072916666666
    4  = 30.983072916666668
    5  = 11.893229166666666
    6  = 9.166666666666666
    7  = 12.619791666666666
    8  = 17.072916666666668
    9  = 15.983072916666666
    10 = 13.346354166666666
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col]
}
